# Apply the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.144.95'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').Value = '3.925.93'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''603.57'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('D6').Value = '''168.03'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.75%  '
$ws.Range('D7').Value = '3.925.15'
$ws.Range('E7').Value = '  +2.22%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').Value = '''0.534'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('E10').Value = '  +1.09%  '
$ws.Range('D11').Value = '''6.44'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').Value = '''0.465'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').Value = '''0.0000256'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.98%  '
$ws.Range('D14').Value = '''37.62'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').Value = '4.582.95'
$ws.Range('E15').Value = '  +2.40%  '
$ws.Range('D16').Value = '3.899.08'
$ws.Range('E16').Value = '  +2.26%  '
$ws.Range('D17').Value = '69.170.93'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('E18').Value = '  -1.28%  '
$ws.Range('D19').Value = '''17.34'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('D21').Value = '''10.96'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -5.74%  '
$ws.Range('D22').Value = '''496.81'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.53%  '
$ws.Range('D23').Value = '''0.731'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').Value = '''0.0000166'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +5.35%  '
$ws.Range('D25').Value = '''85.11'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('D26').Value = '''2.27'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('D27').Value = '''12.14'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.63%  '
$ws.Range('D28').Value = '''10.22'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.41%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = '''2.98'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').Value = '4.077.47'
$ws.Range('E31').Value = '  +2.41%  '
$ws.Range('D32').Value = '''2.39'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.49%  '
$ws.Range('D33').Value = '''7.77'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.18%  '
$ws.Range('D34').Value = '''31.95'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('D35').Value = '3.891.86'
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('D37').Value = '''1.05'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').Value = '''6.00'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.34%  '
$ws.Range('E39').Value = '  -0.84%  '
$ws.Range('D40').Value = '''3.27'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.30%  '
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('D42').Value = '''0.321'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('D43').Value = '''430.96'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.75%  '
$ws.Range('D44').Value = '''2.00'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D45').Value = '''48.08'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.05%  '
$ws.Range('D46').Value = '''8.59'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +2.05%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '''0.000279'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +22.91%  '
$ws.Range('B49').Value = 'Monero'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Value = '''142.97'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '''0.0362'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.87%  '
$ws.Range('D51').Value = '2.797.81'
$ws.Range('E51').Value = '  -1.65%  '
